# Apply the BoM update: two components (D203/D204/D205/D206 LEDs-as-diodes
# and R214/R215/R216/R217 resistors) were removed from the schematic, so the
# corresponding "References" lists and "Quantity Per PCB" counts shrink by 4
# each, and the workbook-wide component totals (Component Count / Fitted
# Components / Total Components) drop by 8 (4+4) on both the BoM and DNF
# sheets.

$wb = $excel.ActiveWorkbook

function Update-ReferenceGroup {
    param(
        [string]$sheetName,
        [string]$removeText,
        [int]$qtyDelta
    )

    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Cells.Find($removeText)

    $refCell = $ws.Cells.Item($cell.Row, $cell.Column)
    $oldRefs = $refCell.Text
    $newRefs = $oldRefs -replace [regex]::Escape(" $removeText"), ""
    $refCell.Value = $newRefs

    $qtyCell = $ws.Cells.Item($cell.Row, 12)
    $oldQty = [int]$qtyCell.Text
    $qtyCell.Value = $oldQty - $qtyDelta
}

# 1. LED/diode group (D102..D401): drop D203 D204 D205 D206 (4 refs, qty 16 -> 12)
Update-ReferenceGroup -sheetName "BoM" -removeText "D203 D204 D205 D206" -qtyDelta 4

# 2. Resistor group (R101..R507): drop R214 R215 R216 R217 (4 refs, qty 42 -> 38)
Update-ReferenceGroup -sheetName "BoM" -removeText "R214 R215 R216 R217" -qtyDelta 4

# 3. Update the workbook-wide summary totals on both the BoM and DNF sheets.
#    8 fitted components (4 + 4) were removed in total.
foreach ($sheetName in @("BoM", "DNF")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $componentCount = $ws.Cells.Item(3, 6)
    $componentCount.Value = [int]$componentCount.Text - 8

    $fittedComponents = $ws.Cells.Item(4, 6)
    $fittedComponents.Value = [int]$fittedComponents.Text - 8

    $totalComponents = $ws.Cells.Item(6, 6)
    $totalComponents.Value = [int]$totalComponents.Text - 8
}
